$d = $word.ActiveDocument

# Locate the paragraph that spells out the id tag across three runs:
#   "<id>"  (Courier New / color 7f6000 / 9pt)
#   "p123r_2"  (plain run, no special rFonts/color/size)
#   "</id>" (Courier New / color 7f6000 / 9pt)
# and collapse them into a single run holding the combined text
# "<id>p123r_2</id>", which inherits the formatting of the "<id>" run
# (the first character of the matched range) exactly like Word does
# when a Find/Replace match spans multiple runs.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*<id>p123r_2</id>*") {
        $p.Range.Find.Execute("<id>p123r_2</id>", $false, $false, $false,
                               $false, $false, $true, 1, $false,
                               "<id>p123r_2</id>", 2)
    }
}
